# Updates Flashscore odds/stats figures in Sheet1 as per the 2025-04-17 data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7
$ws.Range("G7").Value = 2.55
$ws.Range("I7").Value = 2.9
$ws.Range("T7").Value = 6.5
$ws.Range("AG7").Value = 12
$ws.Range("AI7").Value = 29

# Row 12
$ws.Range("L12").Value = 1.33
$ws.Range("M12").Value = 2.77
$ws.Range("P12").Value = 1.42
$ws.Range("Q12").Value = 2.45
$ws.Range("R12").Value = 2.1
$ws.Range("T12").Value = 5.6
$ws.Range("U12").Value = 6.2
$ws.Range("X12").Value = 13.5
$ws.Range("Y12").Value = 35
$ws.Range("AB12").Value = 22
$ws.Range("AJ12").Value = 80

# Row 13
$ws.Range("G13").Value = 2.1
$ws.Range("H13").Value = 3.2
$ws.Range("I13").Value = 3.4
$ws.Range("R13").Value = 1.73
$ws.Range("S13").Value = 2
$ws.Range("U13").Value = 10
$ws.Range("V13").Value = 9
$ws.Range("W13").Value = 19
$ws.Range("X13").Value = 17
$ws.Range("Y13").Value = 29
$ws.Range("AA13").Value = 6
$ws.Range("AB13").Value = 13
$ws.Range("AE13").Value = 10
$ws.Range("AF13").Value = 17
$ws.Range("AG13").Value = 12
$ws.Range("AI13").Value = 29
$ws.Range("AJ13").Value = 34

# Row 14
$ws.Range("G14").Value = 2.32
$ws.Range("I14").Value = 2.77
$ws.Range("P14").Value = 1.4
$ws.Range("Q14").Value = 2.52
$ws.Range("W14").Value = 23
$ws.Range("X14").Value = 19
$ws.Range("AE14").Value = 8.75
$ws.Range("AI14").Value = 24

# Row 16
$ws.Range("G16").Value = 2.05
$ws.Range("I16").Value = 3.5
$ws.Range("T16").Value = 9
$ws.Range("W16").Value = 19
$ws.Range("AH16").Value = 41
$ws.Range("AI16").Value = 26

# Row 17
$ws.Range("G17").Value = 1.4
$ws.Range("H17").Value = 4.75
$ws.Range("I17").Value = 7.5
$ws.Range("N17").Value = 1.67
$ws.Range("O17").Value = 2.15
$ws.Range("P17").Value = 1.27
$ws.Range("Q17").Value = 3.4
$ws.Range("R17").Value = 1.91
$ws.Range("S17").Value = 1.91
$ws.Range("U17").Value = 7
$ws.Range("W17").Value = 9.5
$ws.Range("X17").Value = 11
$ws.Range("Z17").Value = 15
$ws.Range("AA17").Value = 9
$ws.Range("AB17").Value = 19
$ws.Range("AD17").Value = 301
$ws.Range("AE17").Value = 19
$ws.Range("AF17").Value = 41
$ws.Range("AG17").Value = 21
$ws.Range("AH17").Value = 81
$ws.Range("AI17").Value = 51
$ws.Range("AJ17").Value = 51

# Row 19
$ws.Range("I19").Value = 3.4
$ws.Range("J19").Value = 1.02
$ws.Range("K19").Value = 12
$ws.Range("L19").Value = 1.22
$ws.Range("M19").Value = 4
$ws.Range("R19").Value = 1.67
$ws.Range("S19").Value = 2.1
$ws.Range("X19").Value = 15
$ws.Range("Y19").Value = 23
$ws.Range("Z19").Value = 12
$ws.Range("AE19").Value = 12
$ws.Range("AF19").Value = 19
$ws.Range("AG19").Value = 13

# Row 21
$ws.Range("N21").Value = 1.65
$ws.Range("O21").Value = 2.2

